$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values from the refreshed
# symbol-list data pull. Force text format first so Excel stores the
# exact literal strings (matching the source feed formatting, e.g.
# trailing zeros / percent signs) instead of auto-coercing to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '303.08'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '5.48%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '32.00'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '9.66%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.247'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.85%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07456'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '6.81%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.854'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.56%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.793'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '6.80%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.523'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '7.98%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9190'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.97%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01753'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2,590.91%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1686'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.81%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07983'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '5.44%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08007'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.63%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03037'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.15%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09903'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '9.85%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001495'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-4.89%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04618'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2.27%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006213'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.36%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.08%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.232'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.09%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1344'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.85%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.500'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '12.00%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1622'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.45%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001220'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.87%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004448'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4.79%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001401'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '19.75%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001747'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '4.71%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04495'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007178'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.57%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1350'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '8.41%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002211'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '6.84%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01279'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '10.16%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006157'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '5.59%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7097'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-63.22%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01301'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.44%'
